$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated metric values (B, C, D columns) for rows 2-14
$data = @(
    @{ Row = 2;  B = 0.5343502614175931;  C = 0.9893599515376003;  D = 0.6114124767059477 },
    @{ Row = 3;  B = 0.2478145157890329;  C = 0.9951564590568358;  D = 0.3905129923112245 },
    @{ Row = 4;  B = 0.2906691478160842;  C = 0.9944081719399989;  D = 0.4444675901875104 },
    @{ Row = 5;  B = 0.4142025673419836;  C = 0.9918324924658926;  D = 0.4966819928325304 },
    @{ Row = 6;  B = 0.4849316671866268;  C = 0.9857562866251085;  D = 0.5129315018905797 },
    @{ Row = 7;  B = 0.09456831123154207; C = 0.9986930403330508;  D = 0.2513101752851327 },
    @{ Row = 8;  B = 0.03699246652441175; C = 0.9996187265353078;  D = 0.1351892691472226 },
    @{ Row = 9;  B = 0.0968294432200346;  C = 0.9994232146343792;  D = 0.2189325545095534 },
    @{ Row = 10; B = 0.06752566960878822; C = 0.9987664659309182;  D = 0.2059592627517777 },
    @{ Row = 11; B = 0.1273178433545562;  C = 0.9905932939175162;  D = 0.2757376756410942 },
    @{ Row = 12; B = 0.05272969306323207; C = 0.9984493096864386;  D = 0.1681976945500981 },
    @{ Row = 13; B = 0.05907475696115474; C = 0.9994396494559892;  D = 0.1751086826368257 },
    @{ Row = 14; B = 0.05578305852789792; C = 0.9992447318213219;  D = 0.1901039495598666 }
)

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       RandomForestRegressor())]),`n                                            param_grid={'model__max_depth': [3,`n                                                                             5,`n                                                                             7],`n                                                        'model__n_estimators': [50,`n                                                                                100,`n                                                                                150]},`n                                            scoring='neg_mean_squared_error'))"

# Add new header in F1, copying the header formatting from A1 (bold, bordered, centered)
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "Modelo"

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 6).Value = $modelText
}
